$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.056684
$ws.Range("H2").Value = 36.170052
$ws.Range("I2").Value = 0.06307822458376462
$ws.Range("J2").Value = 0.06307822458376462
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 3.498721092170666
$ws.Range("R2").Value = 31.488489829536
$ws.Range("S2").Value = 0.002163198059941946
$ws.Range("T2").Value = 0.002163198059941946
$ws.Range("G3").Value = 12.056684
$ws.Range("H3").Value = 36.170052
$ws.Range("I3").Value = 0.06307822458376462
$ws.Range("J3").Value = 0.06307822458376462
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 85.72725111718933
$ws.Range("R3").Value = 771.5452600547039
$ws.Range("S3").Value = 0.05300366002761503
$ws.Range("T3").Value = 0.05300366002761503
$ws.Range("G4").Value = 12.056684
$ws.Range("H4").Value = 36.170052
$ws.Range("I4").Value = 0.06307822458376462
$ws.Range("J4").Value = 0.06307822458376462
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 12.79571452135867
$ws.Range("R4").Value = 115.161430692228
$ws.Range("S4").Value = 0.007911366496207651
$ws.Range("T4").Value = 0.007911366496207649
$ws.Range("I5").Value = 0.1315309049843414
$ws.Range("J5").Value = 0.1315309049843414
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 7.295543820037334
$ws.Range("R5").Value = 65.659894380336
$ws.Range("S5").Value = 0.004510707147546583
$ws.Range("T5").Value = 0.004510707147546583
$ws.Range("I6").Value = 0.1315309049843414
$ws.Range("J6").Value = 0.1315309049843414
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.1105233924530741
$ws.Range("T6").Value = 0.1105233924530741
$ws.Range("I7").Value = 0.1315309049843414
$ws.Range("J7").Value = 0.1315309049843414
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 26.68166268187534
$ws.Range("S7").Value = 0.01649680538372068
$ws.Range("T7").Value = 0.01649680538372068
$ws.Range("H8").Value = 461.8238680000001
$ws.Range("I8").Value = 0.8053908704318941
$ws.Range("J8").Value = 0.8053908704318941
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 44.67212012411378
$ws.Range("R8").Value = 402.0490811170241
$ws.Range("S8").Value = 0.02761999057376211
$ws.Range("T8").Value = 0.02761999057376211
$ws.Range("H9").Value = 461.8238680000001
$ws.Range("I9").Value = 0.8053908704318941
$ws.Range("J9").Value = 0.8053908704318941
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.6767575366524262
$ws.Range("T9").Value = 0.6767575366524262
$ws.Range("H10").Value = 461.8238680000001
$ws.Range("I10").Value = 0.8053908704318941
$ws.Range("J10").Value = 0.8053908704318941
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 163.3773259180725
$ws.Range("S10").Value = 0.1010133432057058
$ws.Range("T10").Value = 0.1010133432057058
